$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume(1h) (E) columns.
# Price cells are forced to Text (matching the source data, which stores
# numeric-looking prices as text) via NumberFormat "@", then the style is
# reset back to Normal so no stray formatting is introduced.

# Row 2: Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.943.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.00%  "

# Row 3: Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.952.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.97%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.53%  "

# Row 6: Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.19%  "

# Row 7: USDC
$ws.Range("E7").Value = "  +0.04%  "

# Row 8: LidoStakedEther
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.951.24"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.91%  "

# Row 9: XRP
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.505"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.94%  "

# Row 10: Toncoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.15%  "

# Row 11: Dogecoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.150"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.25%  "

# Row 12: Cardano
$ws.Range("E12").Value = "  +0.71%  "

# Row 13: ShibaInu
$ws.Range("E13").Value = "  +5.04%  "

# Row 14: Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.34%  "

# Row 15: TRON
$ws.Range("E15").Value = "  -0.60%  "

# Row 16: WrappedliquidstakedEther2.0
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.442.36"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.03%  "

# Row 17: WrappedBTC
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.976.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.05%  "

# Row 18: Polkadot
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.54%  "

# Row 19: WrappedEther
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.953.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.03%  "

# Row 20: BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "440.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.13%  "

# Row 21: Chainlink
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.48"
$ws.Range("D21").Style = "Normal"

# Row 22: Polygon
$ws.Range("E22").Value = "  -0.59%  "

# Row 23: Uniswap
$ws.Range("E23").Value = "  -0.39%  "

# Row 24: RenderToken
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.75%  "

# Row 25: Litecoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.83"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.94%  "

# Row 26: Fetch.AI
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.28%  "

# Row 27: InternetComputer(DFINITY)
$ws.Range("E27").Value = "  +0.81%  "

# Row 28: Dai
$ws.Range("E28").Value = "  +0.03%  "

# Row 29: ImmutableX
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.57%  "

# Row 30: NEARProtocol
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.37%  "

# Row 31: PancakeSwap
$ws.Range("E31").Value = "  +0.91%  "

# Row 32: PEPE
$ws.Range("E32").Value = "  +18.20%  "

# Row 33: EthereumClassic
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.39"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.67%  "

# Row 34: Hedera
$ws.Range("E34").Value = "  -0.81%  "

# Row 35: FirstDigitalUSD
$ws.Range("E35").Value = "  +0.01%  "

# Row 36: Mantle
$ws.Range("E36").Value = "  -1.68%  "

# Row 37: Filecoin
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.62"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.38%  "

# Row 38: dogwifhat
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.04%  "

# Row 39: OKB
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.73"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.02%  "

# Row 40: Stacks
$ws.Range("E40").Value = "  +2.50%  "

# Row 41: Cosmos
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.31%  "

# Row 42: Kaspa
$ws.Range("E42").Value = "  -3.45%  "

# Row 43: TheGraph
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.279"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.41%  "

# Row 44: Arweave
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.63%  "

# Row 45: Maker
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.708.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.78%  "

# Row 46: Monero
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "135.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.65%  "

# Row 47: VeChain
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0339"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.22%  "

# Row 48: Bittensor
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "359.93"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.41%  "

# Row 49: USDe
$ws.Range("E49").Value = "  -0.01%  "

# Row 50: Stellar
$ws.Range("E50").Value = "  -0.28%  "

# Row 51: InjectiveProtocol
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.04%  "
